$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 69
$ws1.Range("F5").Value = 376
$ws1.Range("F6").Value = 11205
$ws1.Range("F7").Value = 611
$ws1.Range("F8").Value = 100
$ws1.Range("F9").Value = 11
$ws1.Range("F12").Value = 155
$ws1.Range("F15").Value = 44
$ws1.Range("F19").Value = 1217
$ws1.Range("F20").Value = 61

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 69
$ws4.Range("F5").Value = 376
$ws4.Range("F6").Value = 11205
$ws4.Range("F7").Value = 611
$ws4.Range("F8").Value = 100
$ws4.Range("F9").Value = 11
$ws4.Range("F12").Value = 155
$ws4.Range("F15").Value = 44
$ws4.Range("F19").Value = 1217
$ws4.Range("F20").Value = 61
$ws4.Range("F21").Value = 886
